# Automatische test-sync: 2025-07-22 12:37:50
#
# This script reproduces the commit that:
#  1) Adds a new "Testmail #12" row (row 12) to the "Logs" sheet,
#     expanding the used range, conditional formatting ranges and dimension.
#  2) Re-orders / updates the "Dashboard" summary sheet category counts
#     (adds a new "Overig" row, shuffles the existing category order)
#     to stay in sync with the Logs sheet.
#  3) Extends the bar chart on the Dashboard sheet so its category / value
#     series reference the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet - append row 12
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Ik wacht nog steeds op reactie. Wanneer hoor ik iets?"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Testmail #12: Ik wacht nog steeds op reactie. Wanneer hoor ik iets?"
$logs.Range("D12").Value = "Overig"
$logs.Range("E12").Value = "Beste klant,`r`nDank u wel voor uw e-mail. Excuses voor het ongemak dat u heeft ervaren. Om u beter van dienst te kunnen zijn, kunt u alstublieft uw gebruikersnaam vermelden zodat we het specifieke probleem kunnen onderzoeken en een passende oplossing kunnen bieden.`r`nWij streven ernaar om binnen 24 uur op al onze e-mails te reageren. Zodra we meer informatie hebben, zullen we direct contact met u opnemen.`r`nMet vriendelijke groet,`r`n[Naam van het bedrijf] E-mailassistent"
$logs.Range("F12").Value = "2025-07-22 12:37:22"
$logs.Range("G12").Value = "Ja"
$logs.Range("H12").Value = "Nee"
$logs.Range("I12").Value = "Ja"
$logs.Range("J12").Value = "Ja"

# Grow the conditional formatting ranges (D/G/H/I/J) from row 11 to row 12.
$logs.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D12"))
$logs.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G12"))
$logs.Range("H2:H11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H12"))
$logs.Range("I2:I11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I12"))
$logs.Range("J2:J11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J12"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet - refresh category counts (new order + new row)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Retour / Terugbetaling"
$dash.Range("B2").Value = 4

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 4

$dash.Range("A4").Value = "Openingstijden / Locatie"
$dash.Range("B4").Value = 1

$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 1

$dash.Range("A6").Value = "Overig"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3) Chart - extend the category / value series to include the new row
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$series.Values = "='Dashboard'!`$B`$2:`$B`$6"
